# Update the "FROM"/"TO" schedule date (F14/H14) and the hourly generation
# figures (column F, rows 23-46) on each power-station worksheet.
#
# The new shared date/time value used on every sheet (replaces the old
# 44228.15069444444 / 44222.68194444444 stamps):
$newDate = 44305.361805555556

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Edwaleni
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Edwaleni")
$ws.Range("F14").Value = $newDate
$ws.Range("H14").Value = $newDate

$ws.Range("F23").Value = 2.4
$ws.Range("F24").Value = 2.4
$ws.Range("F25").Value = 2.4
$ws.Range("F26").Value = 2.4
$ws.Range("F27").Value = 2.4
$ws.Range("F28").Value = 2.4

$ws.Range("F29").Value = 14.6

$ws.Range("F33").Value = 14.6
$ws.Range("F34").Value = 14.6
$ws.Range("F35").Value = 14.6
$ws.Range("F36").Value = 14.6
$ws.Range("F37").Value = 14.6
$ws.Range("F38").Value = 14.6
$ws.Range("F39").Value = 14.6
$ws.Range("F40").Value = 14.6

$ws.Range("F43").Value = 2.4
$ws.Range("F44").Value = 2.4
$ws.Range("F45").Value = 2.4
$ws.Range("F46").Value = 2.4

# ---------------------------------------------------------------------
# Maguduza
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Maguduza")
$ws.Range("F14").Value = $newDate
$ws.Range("H14").Value = $newDate

$ws.Range("F29").Value = 5.6

$ws.Range("F33").Value = 5.6
$ws.Range("F34").Value = 5.6
$ws.Range("F35").Value = 5.6
$ws.Range("F36").Value = 5.6
$ws.Range("F37").Value = 5.6
$ws.Range("F38").Value = 5.6
$ws.Range("F39").Value = 5.6
$ws.Range("F40").Value = 5.6

# ---------------------------------------------------------------------
# Ezulwini
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Ezulwini")
$ws.Range("F14").Value = $newDate
$ws.Range("H14").Value = $newDate

$ws.Range("F23").Value = 20
$ws.Range("F24").Value = 20
$ws.Range("F25").Value = 20
$ws.Range("F26").Value = 20
$ws.Range("F27").Value = 20
$ws.Range("F28").Value = 20
$ws.Range("F29").Value = 20

$ws.Range("F33").Value = 20
$ws.Range("F34").Value = 20
$ws.Range("F35").Value = 20
$ws.Range("F36").Value = 20
$ws.Range("F37").Value = 20
$ws.Range("F38").Value = 20
$ws.Range("F39").Value = 20
$ws.Range("F40").Value = 20

$ws.Range("F43").Value = 20
$ws.Range("F44").Value = 20
$ws.Range("F45").Value = 20
$ws.Range("F46").Value = 20

# ---------------------------------------------------------------------
# Maguga (F14/H14 previously held formulas referencing Edwaleni; they
# are replaced with the same literal date value used on the other sheets)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Maguga")
$ws.Range("F14").Value = $newDate
$ws.Range("H14").Value = $newDate

$ws.Range("F29").Value = 10

$ws.Range("F30").Value = 20
$ws.Range("F31").Value = 20
$ws.Range("F32").Value = 20

$ws.Range("F41").Value = 20
$ws.Range("F42").Value = 20
